$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update J1/K1: these were text-shared-string cells ("r"/"s"); they
# become plain numeric values, matching the numeric pattern used by the
# rest of the J/K columns (J=0.6, K=0.5).
$ws.Range("J1").Value = 0.6
$ws.Range("K1").Value = 0.5

# --- Update K2:K51 from 0.6 to 0.5 (model retrained -> new threshold column)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = 0.5
}

# --- Update the view/selection state to match the saved workbook:
#     top-left visible cell A40, active selection K52:K54
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("K52:K54").Select()
